$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'318.44"
$ws.Range("E2").Value = "'4.06%"

# Row 3
$ws.Range("D3").Value = "'39.62"
$ws.Range("E3").Value = "'2.06%"

# Row 4
$ws.Range("D4").Value = "'5.139"
$ws.Range("E4").Value = "'0.67%"

# Row 5
$ws.Range("D5").Value = "'0.08210"
$ws.Range("E5").Value = "'2.00%"

# Row 6
$ws.Range("D6").Value = "'2.088"
$ws.Range("E6").Value = "'7.39%"

# Row 7
$ws.Range("D7").Value = "'8.306"
$ws.Range("E7").Value = "'3.87%"

# Row 8
$ws.Range("D8").Value = "'4.296"
$ws.Range("E8").Value = "'2.52%"

# Row 9
$ws.Range("D9").Value = "'0.9329"
$ws.Range("E9").Value = "'0.18%"

# Row 10
$ws.Range("D10").Value = "'0.1395"
$ws.Range("E10").Value = "'-4.13%"

# Row 11
$ws.Range("D11").Value = "'0.1988"
$ws.Range("E11").Value = "'3.64%"

# Row 12
$ws.Range("D12").Value = "'0.09041"
$ws.Range("E12").Value = "'0.41%"

# Row 13
$ws.Range("D13").Value = "'0.03478"
$ws.Range("E13").Value = "'-0.74%"

# Row 14
$ws.Range("D14").Value = "'0.09804"
$ws.Range("E14").Value = "'0.23%"

# Row 15
$ws.Range("D15").Value = "'0.001394"
$ws.Range("E15").Value = "'0.29%"

# Row 16
$ws.Range("D16").Value = "'0.006273"
$ws.Range("E16").Value = "'3.82%"

# Row 17
$ws.Range("D17").Value = "'3.677"
$ws.Range("E17").Value = "'-2.59%"

# Row 18
$ws.Range("D18").Value = "'3.319"
$ws.Range("E18").Value = "'-2.67%"

# Row 19
$ws.Range("D19").Value = "'0.3473"
$ws.Range("E19").Value = "'1.59%"

# Row 20
$ws.Range("D20").Value = "'0.1292"
$ws.Range("E20").Value = "'-3.07%"

# Row 21
$ws.Range("D21").Value = "'4.864"
$ws.Range("E21").Value = "'1.60%"

# Row 23
$ws.Range("D23").Value = "'0.04324"
$ws.Range("E23").Value = "'-1.33%"

# Row 24
$ws.Range("E24").Value = "'-0.95%"

# Row 25
$ws.Range("D25").Value = "'0.004759"
$ws.Range("E25").Value = "'11.34%"

# Row 26
$ws.Range("E26").Value = "'-0.17%"

# Row 27
$ws.Range("D27").Value = "'0.0004000"
$ws.Range("E27").Value = "'-10.06%"

# Row 39
$ws.Range("D39").Value = "'0.02233"
$ws.Range("E39").Value = "'9.30%"

# Row 40
$ws.Range("D40").Value = "'0.05223"
$ws.Range("E40").Value = "'3.90%"

# Row 41
$ws.Range("D41").Value = "'0.007459"
$ws.Range("E41").Value = "'0.44%"

# Row 42
$ws.Range("D42").Value = "'0.009563"
$ws.Range("E42").Value = "'-5.42%"

# Row 43
$ws.Range("E43").Value = "'2.55%"

# Row 44
$ws.Range("D44").Value = "'0.002150"
$ws.Range("E44").Value = "'1.24%"

# Row 45
$ws.Range("D45").Value = "'0.009852"
$ws.Range("E45").Value = "'8.83%"

# Row 46
$ws.Range("D46").Value = "'0.00006603"
$ws.Range("E46").Value = "'6.92%"

# Row 47
$ws.Range("E47").Value = "'-0.07%"

# Row 48
$ws.Range("B48").Value = "'CoinbaseStockToken"
$ws.Range("C48").Value = "'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.001200"
$ws.Range("E48").Value = "'-25.07%"

# Row 49
$ws.Range("B49").Value = "'BOLO"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.002770"
$ws.Range("E49").Value = "'-0.91%"

# Row 50
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.07%"

# Row 51
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.07%"
